$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows scraped for Murali Vijay (CSK) - appended after the existing
# 3 data rows (rows 2-4), re-adding the same 3 matches (rows 5-7).
$newRows = @(
    @{ Row = 5; Values = @(" Abu Dhabi", " September 19 2020", "Super Kings won by 5 wickets (with 4 balls remaining)", "Chennai Super Kings", "Mumbai Indians", "Murali Vijay ", "1", "7", "0", "0", "14.28") },
    @{ Row = 6; Values = @(" Sharjah", " September 22 2020", "Royals won by 16 runs", "Chennai Super Kings", "Rajasthan Royals", "Murali Vijay ", "21", "21", "3", "0", "100.00") },
    @{ Row = 7; Values = @(" Dubai (DSC)", " September 25 2020", "Capitals won by 44 runs", "Chennai Super Kings", "Delhi Capitals", "Murali Vijay ", "10", "15", "1", "0", "66.66") }
)

# Columns A1:K1 -> A..K (1-indexed column numbers 1..11)
foreach ($entry in $newRows) {
    $r = $entry.Row
    $vals = $entry.Values
    for ($col = 1; $col -le $vals.Count; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        $text = $vals[$col - 1]

        # Columns G..K (7..11) hold digit/decimal strings ("1", "21", "14.28", ...).
        # Assigning those through COM Value would auto-coerce to a real number,
        # but the source data keeps them as plain text - force text storage via
        # a Text number format, then strip the format back off so no stray
        # cell style lingers on the saved cell.
        if ($col -ge 7) {
            $cell.NumberFormat = "@"
            $cell.Value = $text
            $cell.ClearFormats()
        } else {
            $cell.Value = $text
        }
    }
}
